# Update children's names in the "Child" column (column A) of the active
# worksheet to their corrected/expanded forms, per the funding extractor
# matching-logic improvements described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "Boden Williams"               = "Boden Nelson Williams"
    "Calum Carroll ward"           = "Calum Thomas Carroll Ward"
    "Charlotte Guyler"             = "Lottie Charlotte Guyler"
    "Charlotte Rose Doyle"         = "Charlotte Doyle"
    "Eleanor Wadden"               = "Eleanor Niamh Wadden"
    "Elsie Williams"               = "Elsie Nelson Williams"
    "Eve Otoole"                   = "Eve O Toole"
    "Jack Vickers McGerr"          = "JP Vickers McGerr"
    "Juno Luna Hynes Byrne"        = "Juno Hynes Byrne"
    "Katie Vickers Mc Gerr"        = "Katie Mcgerr"
    "Lily Grnik"                   = "Lily Gornik"
    "Lily Kathy May Corcoran"      = "Lily Corcoran"
    "Marc Aurele Gaaloul Donnelly" = "Marc Gaaloul Donnelly"
    "Naoise Siochr"                = "Naoise O Siochru"
    "Ray OCleirigh"                = "Ray O Cleirigh"
    "Theo OShaughnessy"            = "Theo O Shaughnessy"
    "Toms Hobbs"                   = "Toms Carmody Finnegan"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
